$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Record the Hyperlink style so we can re-apply it after re-creating hyperlinks
$hyperlinkStyle = $ws.Range("E3").Style

# Record existing hyperlinks (row, column, target address) before the row is removed.
$links = @()
foreach ($h in $ws.Hyperlinks) {
    $links += , @($h.Range.Row, $h.Range.Column, $h.Address)
}

# Select row 7 (Jaxine L Wolfe, metadata Provider) and delete it, shifting rows 8:10 up to 7:9.
$ws.Rows.Item(7).Select()
$ws.Rows.Item(7).Delete()

# The hyperlink collection does not automatically re-target itself after the row shift,
# so rebuild it at the correct (possibly shifted) cells, preserving order/target/style.
$ws.Hyperlinks.Delete()
foreach ($l in $links) {
    $row = $l[0]
    if ($row -gt 7) { $row = $row - 1 }
    $col = $l[1]
    $addr = $l[2]
    $cell = $ws.Cells.Item($row, $col)
    $ws.Hyperlinks.Add($cell, $addr) | Out-Null
    $cell.Style = $hyperlinkStyle
}

# Restore selection to row 7 (now occupied by the former row 8 contents), matching the
# post-edit state in the workbook.
$ws.Rows.Item(7).Select()
